$d = $word.ActiveDocument

# Insert a new paragraph right after the "git commit -m "message":" bullet,
# describing "git commit -a -m "message":" (skips the staging step).
$anchor = $d.Paragraphs.Item(12)
$anchor.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(13)

$r = $newPara.Range
$r.Collapse(1)  # wdCollapseStart

$startPos = $r.Start
$part1 = "git commit –a –m “message”:"
$part2 = " "
$part3 = "It will stage all the tracked modified files and  commit them. But untracked files will not get commit we have to first track them. "

$r.InsertAfter($part1 + $part2 + $part3)

$p1End = $startPos + $part1.Length
$p2End = $p1End + $part2.Length
$p3End = $p2End + $part3.Length

$range1 = $d.Range($startPos, $p1End)
$range1.Font.Bold = $true
$range1.Font.Size = 14

$range2 = $d.Range($p1End, $p2End)
$range2.Font.Bold = $true
$range2.Font.Size = 14

$range3 = $d.Range($p2End, $p3End)
$range3.Font.Size = 14

# Match the paragraph-mark's own run formatting (bold) as recorded in the
# target document, then make sure the descriptive (non-bold) run stays
# un-bolded.
$newPara.Range.Font.Bold = $true
$range3.Font.Bold = $false
